$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing track rows (rows 3 and 4) entirely - shrinks
# the used range down to A1:F2.
$ws.Rows("3:4").Delete()

# Replace row 2's contents with the new track's data.
$ws.Range("B2").Value = '"Dont Panic" lyrics'

$lyrics = "Bones sinking like stones`nAll that we fought for`nHomes, places we've grown`nAll of us are done for`n`nAnd we live in a beautiful world`nYeah we do, yeah we do`nWe live in a beautiful world`n`nBones sinking like stones`nAll that we fought for`nAnd homes, places we've grown`nAll of us are done for`n`nAnd we live in a beautiful world`nYeah we do, yeah we do`nWe live in a beautiful world`n`nAnd we live in a beautiful world`nYeah we do, yeah we do`nWe live in a beautiful world`n`nOh, all that I know`nThere's nothing here to run from`n'Cause yeah everybody here's got somebody to lean on"
$ws.Range("C2").Value = $lyrics

$ws.Range("D2").Value = "`nColdplay Lyrics`n"
$ws.Range("E2").Value = 'EP: "The Blue Room" (1999)'
$ws.Range("F2").Value = "Submit CorrectionsWriter(s): Jonathan Mark Buckland, Christopher Anthony John Martin, William Champion, Guy Rupert Berryman"

# The long multi-line lyrics text makes Excel auto-expand the row height;
# re-fit it back down so row 2 keeps its original (default) height.
$ws.Rows(2).AutoFit()
